# Weekly update: insert a new price record row for "Cilantro" (Agricola del
# Norte S.A. de Arica) as the new first data row after the header block that
# was previously row 78, shifting the remaining historical rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 78; this shifts old rows 78-112 down to 79-113
# and copies formatting (including the date style on column D) from row 78.
$ws.Rows.Item(78).Insert()

# Populate the newly inserted row 78 with the new weekly record. All of the
# descriptive columns (A, B, C, E, F, G, H, I, J, N, O, Q, R) are identical to
# every other record in this sub-table, only the date (D) and the price
# columns (K, L, M, P) differ for this new entry.
$row = 78
$ws.Cells.Item($row, 1).Value = 1
$ws.Cells.Item($row, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item($row, 3).Value = "Arica y Parinacota"
$ws.Cells.Item($row, 4).Value = 45027
$ws.Cells.Item($row, 5).Value = 15
$ws.Cells.Item($row, 6).Value = 100112040
$ws.Cells.Item($row, 7).Value = "Cilantro"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 300
$ws.Cells.Item($row, 11).Value = 2000
$ws.Cells.Item($row, 12).Value = 2500
$ws.Cells.Item($row, 13).Value = 2250
$ws.Cells.Item($row, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item($row, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($row, 16).Value = 1125
$ws.Cells.Item($row, 17).Value = 2
$ws.Cells.Item($row, 18).Value = "Hortaliza"
